# Auto-generated PowerShell COM-interop script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update author/title cell (A2): Albert Einstein -> Francois Burgener ---
$ws.Range("A2").Value = "François Burgener"

# --- Update journal rows (dates, activities, hours) ---
# Insertion order chosen to reproduce the shared-strings table order of the target file.
$ws.Range("B8").Value = "Discution du choix de l'architecture de notre application"
$ws.Range("B5").Value = "Discussion sur la proposition du sujets"
$ws.Range("B6").Value = "Lecture du sujet"
$ws.Range("B9").Value = "Rédaction détaillé des fonctionnalités (Budget,catégorisation et choix de devises) et discution sur les necessité du projet"
$ws.Range("B11").Value = "Elaboration du cahier des charges, discussion en groupe pour l'achèvement des fonctionnalités et le choix de l'architecture"
$ws.Range("B7").Value = "Feedback et discution avec le professeur de notre sujet"
$ws.Range("B13").Value = "Finalisation du planning et la distribution des heure. Relecture du cahier des charges"
$ws.Range("B12").Value = "Elaboration des croquis des fenêtes des interface graphique  "
$ws.Range("B16").Value = "Discution du cahier des charges avec le professeur"
$ws.Range("B10").Value = "Elaboration des croquis des fenêtes des interface graphique (vue global,Budget,vue budget,transaction)  "
$ws.Range("B10").Characters(104,1).Font.Color = 255
$ws.Range("B17").Value = "Discution,correction sur les maquettes ainsi que notre schéma relationel"
$ws.Range("B15").Value = "Apprentissage javaFX. Visionnage de tutoriel sur youtube et quelque test "
$ws.Range("B14").Value = "Apprentissage javaFX. J'ai lu le cours de openclassroom sur javaFX. J'ai regardé des vidéo et j'ai fait des petit test pour prendre en main  le code"
$ws.Range("B18").Value = "Apprentissage javaFX, interface avec FXML au lieu du code. Début de la fenêtre de connexion et création d'un compte utilisateur"
$ws.Range("B19").Value = "Continuer sur la fenêtre de connexion, pkus sur le design "
$ws.Range("B20").Value = "Quelque correction a propos des containers utilisés dans la vue de login/register"

# --- Set dates, hours, and row heights for each row ---
$ws.Range("A5").Value = 43150
$ws.Range("C5").Value = 0.5
$ws.Range("A6").Value = 43152
$ws.Range("C6").Value = 0.25
$ws.Range("A7").Value = 43157
$ws.Range("C7").Value = 0.5
$ws.Range("A8").Value = 43157
$ws.Range("C8").Value = 1
$ws.Range("A9").Value = 43158
$ws.Range("C9").Value = 2
$ws.Rows(9).RowHeight = 45
$ws.Range("A10").Value = 43161
$ws.Range("C10").Value = 3
$ws.Rows(10).RowHeight = 30
$ws.Range("A11").Value = 43163
$ws.Range("C11").Value = 4.5
$ws.Rows(11).RowHeight = 45
$ws.Range("A12").Value = 43163
$ws.Range("C12").Value = 2
$ws.Rows(12).RowHeight = 30
$ws.Range("A13").Value = 43164
$ws.Range("C13").Value = 1.5
$ws.Rows(13).RowHeight = 30
$ws.Range("A14").Value = 43169
$ws.Range("C14").Value = 2.5
$ws.Rows(14).RowHeight = 45
$ws.Range("A15").Value = 43170
$ws.Range("C15").Value = 2
$ws.Rows(15).RowHeight = 30
$ws.Range("A16").Value = 43171
$ws.Range("C16").Value = 0.5
$ws.Range("A17").Value = 43171
$ws.Range("C17").Value = 1
$ws.Rows(17).RowHeight = 30
$ws.Range("A18").Value = 43176
$ws.Range("C18").Value = 5
$ws.Rows(18).RowHeight = 45
$ws.Range("A19").Value = 43177
$ws.Range("C19").Value = 1.5
$ws.Range("A20").Value = 43178
$ws.Range("C20").Value = 0.5
$ws.Rows(20).RowHeight = 30

# --- Recalculate total (formula already in C32, engine recalcs automatically) ---
$wb.Application.Calculate()

# --- Update selection to match target view state ---
[void]$ws.Range("G21").Select()

Write-Host "Edit complete"
